$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 5291.0684
$ws.Range("I15").Value = 5291.0684
$ws.Range("K15").Value = 15873.2052
$ws.Range("M15").Value = -15704.2052
$ws.Range("H17").Value = 308
$ws.Range("J17").Value = 334.8846
$ws.Range("L17").Value = 1004.6538
$ws.Range("N17").Value = -1340.6538
$ws.Range("H18").Value = 62501548
$ws.Range("I18").Value = 100001280
$ws.Range("K18").Value = 100001280
$ws.Range("M18").Value = -100000996
$ws.Range("H19").Value = 7541.3335
$ws.Range("J19").Value = 7437.25
$ws.Range("L19").Value = 7437.25
$ws.Range("N19").Value = -7787.25
$ws.Range("H40").Value = 4233.1665
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4825
$ws.Range("H51").Value = 4285.7144
$ws.Range("I51").Value = 4000
$ws.Range("J51").Value = 4500
$ws.Range("K51").Value = 4000
$ws.Range("L51").Value = 4500
$ws.Range("M51").Value = -3516
$ws.Range("N51").Value = -5468
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H112").Value = 1147.2222
$ws.Range("J112").Value = 1184.375
$ws.Range("L112").Value = 3553.125
$ws.Range("N112").Value = -5769.125
$ws.Range("H132").Value = 5178.467
$ws.Range("I132").Value = 3386.318
$ws.Range("J132").Value = 10106.875
$ws.Range("K132").Value = 10158.954
$ws.Range("L132").Value = 30320.625
$ws.Range("M132").Value = -7628.954000000002
$ws.Range("N132").Value = -35380.625

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 5000
$ws.Range("I8").Value = 5000
$ws.Range("J8").Value = 5000
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 5000
$ws.Range("M8").Value = -4856
$ws.Range("N8").Value = -5288
$ws.Range("H63").Value = 2311.3684
$ws.Range("I63").Value = 1546.2222
$ws.Range("K63").Value = 1546.2222
$ws.Range("M63").Value = -860.2221999999999
$ws.Range("H66").Value = 2311.3684
$ws.Range("I66").Value = 1546.2222
$ws.Range("K66").Value = 7731.111
$ws.Range("M66").Value = -4299.111
$ws.Range("H68").Value = 49999
$ws.Range("J68").Value = 49999
$ws.Range("L68").Value = 49999
$ws.Range("N68").Value = -51621
$ws.Range("H71").Value = 49999
$ws.Range("J71").Value = 49999
$ws.Range("L71").Value = 149997
$ws.Range("N71").Value = -158109
$ws.Range("H93").Value = 34969.5
$ws.Range("J93").Value = 34969.5
$ws.Range("L93").Value = 34969.5
$ws.Range("N93").Value = -39961.5
$ws.Range("H102").Value = 6553.2
$ws.Range("J102").Value = 8131.375
$ws.Range("L102").Value = 8131.375
$ws.Range("N102").Value = -11375.375

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 24176.445
$ws.Range("I82").Value = 13942.714
$ws.Range("J82").Value = 59994.5
$ws.Range("K82").Value = 13942.714
$ws.Range("L82").Value = 59994.5
$ws.Range("M82").Value = -13559.714
$ws.Range("N82").Value = -60760.5
$ws.Range("H85").Value = 24176.445
$ws.Range("I85").Value = 13942.714
$ws.Range("J85").Value = 59994.5
$ws.Range("K85").Value = 13942.714
$ws.Range("L85").Value = 59994.5
$ws.Range("M85").Value = -12616.714
$ws.Range("N85").Value = -62646.5
$ws.Range("H86").Value = 1846.2727
$ws.Range("I86").Value = 1773.4375
$ws.Range("K86").Value = 1773.4375
$ws.Range("M86").Value = -650.4375
$ws.Range("H88").Value = 21249.75
$ws.Range("J88").Value = 21249.75
$ws.Range("L88").Value = 21249.75
$ws.Range("N88").Value = -22061.75
$ws.Range("H89").Value = 1846.2727
$ws.Range("I89").Value = 1773.4375
$ws.Range("K89").Value = 8867.1875
$ws.Range("M89").Value = -3251.1875
$ws.Range("H91").Value = 21249.75
$ws.Range("J91").Value = 21249.75
$ws.Range("L91").Value = 21249.75
$ws.Range("N91").Value = -24057.75
$ws.Range("H105").Value = 3962.963
$ws.Range("I105").Value = 2354.1667
$ws.Range("J105").Value = 7180.5557
$ws.Range("K105").Value = 2354.1667
$ws.Range("L105").Value = 7180.5557
$ws.Range("M105").Value = -607.1667000000002
$ws.Range("N105").Value = -10674.5557
$ws.Range("H134").Value = 2225.9321
$ws.Range("I134").Value = 1809.9423
$ws.Range("J134").Value = 5316.143
$ws.Range("K134").Value = 5429.8269
$ws.Range("L134").Value = 15948.429
$ws.Range("M134").Value = -2894.8269
$ws.Range("N134").Value = -21018.429

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4606.6177
$ws.Range("I4").Value = 2824.875
$ws.Range("J4").Value = 5154.846
$ws.Range("K4").Value = 2824.875
$ws.Range("L4").Value = 5154.846
$ws.Range("M4").Value = -2712.875
$ws.Range("N4").Value = -5378.846
$ws.Range("H99").Value = 11737.55
$ws.Range("I99").Value = 9175.857
$ws.Range("J99").Value = 14568.895
$ws.Range("K99").Value = 9175.857
$ws.Range("L99").Value = 14568.895
$ws.Range("M99").Value = -7677.857
$ws.Range("N99").Value = -17564.895
$ws.Range("H105").Value = 2054.8235
$ws.Range("I105").Value = 2110.077
$ws.Range("J105").Value = 1875.25
$ws.Range("K105").Value = 2110.077
$ws.Range("L105").Value = 1875.25
$ws.Range("M105").Value = -363.0770000000002
$ws.Range("N105").Value = -5369.25
$ws.Range("H126").Value = 11737.55
$ws.Range("I126").Value = 9175.857
$ws.Range("J126").Value = 14568.895
$ws.Range("K126").Value = 27527.571
$ws.Range("L126").Value = 43706.685
$ws.Range("M126").Value = -25057.571
$ws.Range("N126").Value = -48646.685
$ws.Range("H132").Value = 8608.281999999999
$ws.Range("J132").Value = 2426.75
$ws.Range("L132").Value = 7280.25
$ws.Range("N132").Value = -12340.25
$ws.Range("H134").Value = 4163.2173
$ws.Range("I134").Value = 3290.9756
$ws.Range("J134").Value = 11315.6
$ws.Range("K134").Value = 9872.926800000001
$ws.Range("L134").Value = 33946.8
$ws.Range("M134").Value = -7337.926800000001
$ws.Range("N134").Value = -39016.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1774.8
$ws.Range("I3").Value = 1774.8
$ws.Range("K3").Value = 5324.4
$ws.Range("M3").Value = -5212.4
$ws.Range("H23").Value = 2618.8
$ws.Range("J23").Value = 5348
$ws.Range("L23").Value = 16044
$ws.Range("N23").Value = -16514
$ws.Range("H118").Value = 435.14285
$ws.Range("I118").Value = 435.14285
$ws.Range("K118").Value = 1305.42855
$ws.Range("M118").Value = -62.42855000000009
$ws.Range("H119").Value = 5957.5
$ws.Range("I119").Value = 2966
$ws.Range("K119").Value = 8898
$ws.Range("M119").Value = -4060
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("H141").Value = 1255.7142
$ws.Range("I141").Value = 1255.7142
$ws.Range("K141").Value = 3767.1426
$ws.Range("M141").Value = 1412.8574

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 34375
$ws.Range("I41").Value = 34375
$ws.Range("K41").Value = 34375
$ws.Range("M41").Value = -34020
$ws.Range("H132").Value = 2990.1155
$ws.Range("I132").Value = 2446.9285
$ws.Range("J132").Value = 5271.5
$ws.Range("K132").Value = 7340.7855
$ws.Range("L132").Value = 15814.5
$ws.Range("M132").Value = -4810.7855
$ws.Range("N132").Value = -20874.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2519.5454
$ws.Range("J46").Value = 2539.3333
$ws.Range("L46").Value = 2539.3333
$ws.Range("N46").Value = -2915.3333
$ws.Range("H55").Value = 721.21875
$ws.Range("I55").Value = 504.1875
$ws.Range("J55").Value = 938.25
$ws.Range("K55").Value = 504.1875
$ws.Range("L55").Value = 938.25
$ws.Range("M55").Value = -331.1875
$ws.Range("N55").Value = -1284.25
$ws.Range("H61").Value = 7407.5625
$ws.Range("I61").Value = 11182.2
$ws.Range("J61").Value = 1116.5
$ws.Range("K61").Value = 11182.2
$ws.Range("L61").Value = 1116.5
$ws.Range("M61").Value = -10980.2
$ws.Range("N61").Value = -1520.5
$ws.Range("H68").Value = 1984
$ws.Range("I68").Value = 2046.8235
$ws.Range("K68").Value = 2046.8235
$ws.Range("M68").Value = -1297.8235
$ws.Range("H71").Value = 1984
$ws.Range("I71").Value = 2046.8235
$ws.Range("K71").Value = 10234.1175
$ws.Range("M71").Value = -6490.1175
$ws.Range("H93").Value = 4182.5
$ws.Range("I93").Value = 4182.5
$ws.Range("K93").Value = 4182.5
$ws.Range("M93").Value = -2934.5
$ws.Range("H113").Value = 7407.5625
$ws.Range("I113").Value = 11182.2
$ws.Range("J113").Value = 1116.5
$ws.Range("K113").Value = 11182.2
$ws.Range("L113").Value = 1116.5
$ws.Range("M113").Value = -9012.200000000001
$ws.Range("N113").Value = -5456.5
$ws.Range("H132").Value = 46941.08
$ws.Range("I132").Value = 54987.617
$ws.Range("K132").Value = 164962.851
$ws.Range("M132").Value = -162432.851
